$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.432.80"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.874.41"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7146"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.13"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07892"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3091"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08253"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7261"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.260"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "1.857.31"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.15"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "29.434.23"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.882"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "245.08"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007828"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "2.122.19"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.064"
$ws.Range("E22").Value = "  +6.64%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1619"
$ws.Range("E25").Value = "  +13.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.48"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.027"
$ws.Range("E27").Value = "  +1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.29"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.404"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.100"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05198"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.942"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.194"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7249"
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.673"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01864"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.695"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.187.38"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9067"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.140"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.79"
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.18"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "2.020.81"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5289"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.792"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.907"
$ws.Range("E49").Value = "  +7.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.305"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4310"
$ws.Range("E51").Value = "  +1.19%  "
